# Auto-generated script applying scheduled-runner price updates to the leve profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 42.714287
$ws.Range("I5").Value = 42.714287
$ws.Range("K5").Value = 42.714287
$ws.Range("M5").Value = 72.285713
$ws.Range("H32").Value = 863.63635
$ws.Range("J32").Value = 1083.3334
$ws.Range("L32").Value = 1083.3334
$ws.Range("N32").Value = -1735.3334
$ws.Range("H34").Value = 1499.5
$ws.Range("I34").Value = 1499.5
$ws.Range("K34").Value = 1499.5
$ws.Range("M34").Value = -1296.5
$ws.Range("H36").Value = 1499.5
$ws.Range("I36").Value = 1499.5
$ws.Range("K36").Value = 1499.5
$ws.Range("M36").Value = -784.5
$ws.Range("H40").Value = 4174.4375
$ws.Range("J40").Value = 4424.4165
$ws.Range("L40").Value = 4424.4165
$ws.Range("N40").Value = -4774.4165
$ws.Range("H112").Value = 2138.8
$ws.Range("J112").Value = 2284.1428
$ws.Range("L112").Value = 6852.428400000001
$ws.Range("N112").Value = -9068.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2510.4666
$ws.Range("I2").Value = 1764.4
$ws.Range("K2").Value = 1764.4
$ws.Range("M2").Value = -1651.4
$ws.Range("H63").Value = 3014.4546
$ws.Range("J63").Value = 3900
$ws.Range("L63").Value = 3900
$ws.Range("N63").Value = -5272
$ws.Range("H66").Value = 3014.4546
$ws.Range("J66").Value = 3900
$ws.Range("L66").Value = 19500
$ws.Range("N66").Value = -26364
$ws.Range("H116").Value = 2510.4666
$ws.Range("I116").Value = 1764.4
$ws.Range("K116").Value = 1764.4
$ws.Range("M116").Value = 529.5999999999999
$ws.Range("H132").Value = 1143.2
$ws.Range("I132").Value = 1143.2
$ws.Range("K132").Value = 3429.6
$ws.Range("M132").Value = -899.6000000000004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2510.4666
$ws.Range("I3").Value = 1764.4
$ws.Range("K3").Value = 1764.4
$ws.Range("M3").Value = -1650.4
$ws.Range("H20").Value = 8747.75
$ws.Range("I20").Value = 9164.333000000001
$ws.Range("K20").Value = 9164.333000000001
$ws.Range("M20").Value = -8917.333000000001
$ws.Range("H22").Value = 304.18182
$ws.Range("I22").Value = 282.77777
$ws.Range("K22").Value = 282.77777
$ws.Range("M22").Value = -109.77777
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("H54").Value = 2499
$ws.Range("I54").Value = 1998
$ws.Range("J54").Value = 3000
$ws.Range("K54").Value = 1998
$ws.Range("L54").Value = 3000
$ws.Range("M54").Value = -1514
$ws.Range("N54").Value = -3968
$ws.Range("N31").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 17988.889
$ws.Range("I41").Value = 13500
$ws.Range("J41").Value = 19271.428
$ws.Range("K41").Value = 13500
$ws.Range("L41").Value = 19271.428
$ws.Range("M41").Value = -13072
$ws.Range("N41").Value = -20127.428
$ws.Range("H59").Value = 28141.273
$ws.Range("I59").Value = 22425.666
$ws.Range("J59").Value = 35000
$ws.Range("K59").Value = 22425.666
$ws.Range("L59").Value = 35000
$ws.Range("M59").Value = -21280.666
$ws.Range("N59").Value = -37290
$ws.Range("H99").Value = 7966.6665
$ws.Range("I99").Value = 8900
$ws.Range("K99").Value = 8900
$ws.Range("M99").Value = -7402
$ws.Range("H126").Value = 7966.6665
$ws.Range("I126").Value = 8900
$ws.Range("K126").Value = 26700
$ws.Range("M126").Value = -24230
$ws.Range("H132").Value = 1256.1538
$ws.Range("I132").Value = 1256.1538
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3768.4614
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1238.4614
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6127.625
$ws.Range("I137").Value = 1496.6666
$ws.Range("J137").Value = 8906.200000000001
$ws.Range("K137").Value = 4489.9998
$ws.Range("L137").Value = 26718.6
$ws.Range("M137").Value = 610.0002000000004
$ws.Range("N137").Value = -36918.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2881.75
$ws.Range("I102").Value = 2881.75
$ws.Range("K102").Value = 2881.75
$ws.Range("M102").Value = -1259.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 1000
$ws.Range("K30").Value = 1000
$ws.Range("M30").Value = -892
$ws.Range("H35").Value = 1883.75
$ws.Range("I35").Value = 845
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 845
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = -509
$ws.Range("N35").Value = -5672
$ws.Range("H46").Value = 3088.8572
$ws.Range("J46").Value = 3249.4546
$ws.Range("L46").Value = 3249.4546
$ws.Range("N46").Value = -3625.4546
$ws.Range("H55").Value = 1346
$ws.Range("J55").Value = 493.33334
$ws.Range("L55").Value = 493.33334
$ws.Range("N55").Value = -839.33334
$ws.Range("H61").Value = 3346.2856
$ws.Range("I61").Value = 3295.6667
$ws.Range("K61").Value = 3295.6667
$ws.Range("M61").Value = -3093.6667
$ws.Range("H113").Value = 3346.2856
$ws.Range("I113").Value = 3295.6667
$ws.Range("K113").Value = 3295.6667
$ws.Range("M113").Value = -1125.6667
$ws.Range("H122").Value = 3615.1428
$ws.Range("I122").Value = 3306
$ws.Range("K122").Value = 9918
$ws.Range("M122").Value = -7468

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6624.75
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 8166.3335
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 16332.667
$ws.Range("M81").Value = -2939
$ws.Range("N81").Value = -18454.667
$ws.Range("H84").Value = 6624.75
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 8166.3335
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 81663.33499999999
$ws.Range("M84").Value = -14696
$ws.Range("N84").Value = -92271.33499999999
$ws.Range("H132").Value = 2491.9092
$ws.Range("I132").Value = 2491.9092
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7475.7276
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4945.7276
$ws.Range("N132").ClearContents()
